# Master data and user info
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4: Testuser -> siddesh.mainkar ---
$ws.Cells.Item(4,2).Value = "siddesh.mainkar"
$ws.Cells.Item(4,3).Value = "Siddesh@123"
$ws.Cells.Item(4,4).Value = "Asst. Executive QC"
$ws.Cells.Item(4,5).Value = 9545299529
$ws.Cells.Item(4,6).Value = "siddesh@swanson.co.in"
$ws.Cells.Item(4,7).Value = "SPI0108"

# --- Row 5: Nanliu1 -> pravin.parab ---
$ws.Cells.Item(5,2).Value = "pravin.parab"
$ws.Cells.Item(5,3).Value = "Pravin@123"
$ws.Cells.Item(5,4).Value = "QC Inspector"
$ws.Cells.Item(5,5).Value = 9834706069
$ws.Cells.Item(5,6).Value = "pravin@swanson.co.in"
$ws.Cells.Item(5,7).Value = "SPI0173"

# --- Row 6: Nanliu2 -> vasudevan.chari ---
$ws.Cells.Item(6,2).Value = "vasudevan.chari"
$ws.Cells.Item(6,3).Value = "Vasudevan@123"
$ws.Cells.Item(6,4).Value = "QA Officer"
$ws.Cells.Item(6,5).Value = 9764658004
$ws.Cells.Item(6,6).Value = "qa@swanson.co.in"
$ws.Cells.Item(6,7).Value = "SPI0140"

# --- Row 7: Nanliu3 -> suvarna.parab ---
$ws.Cells.Item(7,2).Value = "suvarna.parab"
$ws.Cells.Item(7,3).Value = "Suvarna@123"
$ws.Cells.Item(7,4).Value = "QC Executive"
$ws.Cells.Item(7,5).Value = 8888884754
$ws.Cells.Item(7,6).Value = "suvarna@swanson.co.in"
$ws.Cells.Item(7,7).Value = "SPI0141"

# Remove the old mailto hyperlinks on F5:F7 (targets have changed / no longer linked)
$ws.Hyperlinks.Delete()

# --- New rows 8-11 ---
$ws.Cells.Item(8,1).Value = 7
$ws.Cells.Item(8,2).Value = "asmita.tari"
$ws.Cells.Item(8,3).Value = "Asmita@123"
$ws.Cells.Item(8,4).Value = "Asst. Executive QC"
$ws.Cells.Item(8,5).Value = 8459474919
$ws.Cells.Item(8,6).Value = "asmita@swanson.co.in"
$ws.Cells.Item(8,7).Value = "SPI0098"

$ws.Cells.Item(9,1).Value = 8
$ws.Cells.Item(9,2).Value = "babaji.gaonkar"
$ws.Cells.Item(9,3).Value = "Babaji@123"
$ws.Cells.Item(9,4).Value = "Asst. Manager QA"
$ws.Cells.Item(9,5).Value = 7798687267
$ws.Cells.Item(9,6).Value = "babaji@swanson.co.in"
$ws.Cells.Item(9,7).Value = "SPI0181"

$ws.Cells.Item(10,1).Value = 9
$ws.Cells.Item(10,2).Value = "raghu.baykar"
$ws.Cells.Item(10,3).Value = "Raghu@123"
$ws.Cells.Item(10,4).Value = "QC Inspector"
$ws.Cells.Item(10,5).Value = 9370230501
$ws.Cells.Item(10,6).Value = "spiqc@swanson.co.in"
$ws.Cells.Item(10,7).Value = "SPI0185"

$ws.Cells.Item(11,1).Value = 10
$ws.Cells.Item(11,2).Value = "manohar.satkhalkar"
$ws.Cells.Item(11,3).Value = "Manohar@123"
$ws.Cells.Item(11,4).Value = "Supervisor"
$ws.Cells.Item(11,5).Value = 8805018631
$ws.Cells.Item(11,6).Value = "manohar@swanson.co.in"
$ws.Cells.Item(11,7).Value = "SPI0124"

# EmployeeID column (G) is left-aligned throughout (header + all data rows)
$ws.Range("G1:G11").HorizontalAlignment = -4131

Write-Host "Data rows updated"
